# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "conversion" note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.99 = 11152.58 pesos`n✅ 11152.58 pesos = 3.0 = 969.62 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate table values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 334.9
$ws2.Range("O10").Value = 3735
$ws2.Range("N12").Value = 3714
$ws2.Range("O12").Value = 322.9
